$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Walk the rows from the bottom up and delete the rows whose first
# cell text starts with "Technical Development" or "Compilance",
# since removing a row shifts the indices of subsequent rows.
for ($i = $t.Rows.Count; $i -ge 1; $i--) {
    $row = $t.Rows.Item($i)
    $label = $row.Cells.Item(1).Range.Text
    if ($label -match "Technical Development" -or $label -match "Compilance") {
        $row.Delete()
    }
}
